# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) column F, and flip a couple of
# sold-out ("已售罄") G-column placeholders over to their now-known ticket
# price, across the 展览 / 演出 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 518
$ws.Cells.Item(6, 6).Value = 934
$ws.Cells.Item(8, 6).Value = 11
$ws.Cells.Item(9, 6).Value = 971
$ws.Cells.Item(10, 6).Value = 760
$ws.Cells.Item(11, 6).Value = 211
$ws.Cells.Item(14, 6).Value = 794
$ws.Cells.Item(15, 6).Value = 261
$ws.Cells.Item(16, 6).Value = 563
$ws.Cells.Item(17, 6).Value = 495
$ws.Cells.Item(18, 6).Value = 1309
$ws.Cells.Item(21, 6).Value = 1132
$ws.Cells.Item(22, 6).Value = 2819
$ws.Cells.Item(23, 6).Value = 1338
$ws.Cells.Item(24, 6).Value = 669
$ws.Cells.Item(26, 6).Value = 1253
$ws.Cells.Item(28, 6).Value = 983
$ws.Cells.Item(29, 6).Value = 330
$ws.Cells.Item(30, 6).Value = 1902
$ws.Cells.Item(31, 6).Value = 38
$ws.Cells.Item(32, 6).Value = 6
$ws.Cells.Item(33, 6).Value = 1351

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 73
$ws.Cells.Item(3, 7).Value = 520
$ws.Cells.Item(4, 6).Value = 353
$ws.Cells.Item(5, 6).Value = 9

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 724

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 724
$ws.Cells.Item(4, 6).Value = 518
$ws.Cells.Item(5, 6).Value = 73
$ws.Cells.Item(8, 7).Value = 520
$ws.Cells.Item(9, 7).Value = 520
$ws.Cells.Item(10, 6).Value = 353
$ws.Cells.Item(11, 6).Value = 9
$ws.Cells.Item(13, 6).Value = 934
$ws.Cells.Item(16, 6).Value = 11
$ws.Cells.Item(17, 6).Value = 971
$ws.Cells.Item(18, 6).Value = 760
$ws.Cells.Item(19, 6).Value = 211
$ws.Cells.Item(27, 6).Value = 794
$ws.Cells.Item(28, 6).Value = 261
$ws.Cells.Item(29, 6).Value = 563
$ws.Cells.Item(30, 6).Value = 495
$ws.Cells.Item(31, 6).Value = 1309
$ws.Cells.Item(34, 6).Value = 1132
$ws.Cells.Item(35, 6).Value = 2819
$ws.Cells.Item(36, 6).Value = 1338
$ws.Cells.Item(37, 6).Value = 669
$ws.Cells.Item(39, 6).Value = 1253
$ws.Cells.Item(43, 6).Value = 983
$ws.Cells.Item(44, 6).Value = 330
$ws.Cells.Item(45, 6).Value = 1902
$ws.Cells.Item(46, 6).Value = 38
$ws.Cells.Item(47, 6).Value = 6
$ws.Cells.Item(48, 6).Value = 1351
